$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1064305.9
$ws.Range("J17").Value = 1064305.9
$ws.Range("L17").Value = 3192917.7
$ws.Range("N17").Value = -3193253.7
$ws.Range("H32").Value = 1242.8572
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 1242.8572
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 1242.8572
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -1894.8572
$ws.Range("H117").Value = 35900
$ws.Range("J117").Value = 35900
$ws.Range("L117").Value = 35900
$ws.Range("N117").Value = -45078
$ws.Range("H137").Value = 1522.7812
$ws.Range("I137").Value = 1065.1072
$ws.Range("J137").Value = 4726.5
$ws.Range("K137").Value = 3195.3216
$ws.Range("L137").Value = 14179.5
$ws.Range("M137").Value = -645.3215999999998
$ws.Range("N137").Value = -19279.5
$ws.Range("H138").Value = 1876.6129
$ws.Range("I138").Value = 1331.5333
$ws.Range("J138").Value = 3319.4707
$ws.Range("K138").Value = 3994.5999
$ws.Range("L138").Value = 9958.4121
$ws.Range("M138").Value = 1145.4001
$ws.Range("N138").Value = -20238.4121

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1509.65
$ws.Range("I2").Value = 1701.0625
$ws.Range("J2").Value = 744
$ws.Range("K2").Value = 1701.0625
$ws.Range("L2").Value = 744
$ws.Range("M2").Value = -1588.0625
$ws.Range("N2").Value = -970
$ws.Range("H4").Value = 146.66667
$ws.Range("I4").Value = 152
$ws.Range("J4").Value = 120
$ws.Range("K4").Value = 152
$ws.Range("L4").Value = 120
$ws.Range("M4").Value = -36
$ws.Range("N4").Value = -352
$ws.Range("H45").Value = 1669.0667
$ws.Range("I45").Value = 1448.4445
$ws.Range("K45").Value = 1448.4445
$ws.Range("M45").Value = -1071.4445
$ws.Range("H97").Value = 967.52
$ws.Range("I97").Value = 694.875
$ws.Range("K97").Value = 694.875
$ws.Range("M97").Value = -198.875
$ws.Range("H110").Value = 24065.5
$ws.Range("I110").Value = 26023.05
$ws.Range("J110").Value = 4490
$ws.Range("K110").Value = 26023.05
$ws.Range("L110").Value = 4490
$ws.Range("M110").Value = -23978.05
$ws.Range("N110").Value = -8580
$ws.Range("H116").Value = 1509.65
$ws.Range("I116").Value = 1701.0625
$ws.Range("J116").Value = 744
$ws.Range("K116").Value = 1701.0625
$ws.Range("L116").Value = 744
$ws.Range("M116").Value = 592.9375
$ws.Range("N116").Value = -5332

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1509.65
$ws.Range("I3").Value = 1701.0625
$ws.Range("J3").Value = 744
$ws.Range("K3").Value = 1701.0625
$ws.Range("L3").Value = 744
$ws.Range("M3").Value = -1587.0625
$ws.Range("N3").Value = -972
$ws.Range("H20").Value = 1889.6
$ws.Range("I20").Value = 1981.3334
$ws.Range("J20").Value = 1752
$ws.Range("K20").Value = 1981.3334
$ws.Range("L20").Value = 1752
$ws.Range("M20").Value = -1734.3334
$ws.Range("N20").Value = -2246
$ws.Range("H35").Value = 38985.715
$ws.Range("I35").Value = 40500
$ws.Range("J35").Value = 38733.332
$ws.Range("K35").Value = 40500
$ws.Range("L35").Value = 38733.332
$ws.Range("M35").Value = -40190
$ws.Range("N35").Value = -39353.332
$ws.Range("H94").Value = 781.4167
$ws.Range("I94").Value = 532.8333
$ws.Range("J94").Value = 1030
$ws.Range("K94").Value = 532.8333
$ws.Range("L94").Value = 1030
$ws.Range("M94").Value = -81.83330000000001
$ws.Range("N94").Value = -1932
$ws.Range("H99").Value = 1358.5476
$ws.Range("I99").Value = 928.44446
$ws.Range("J99").Value = 2132.7334
$ws.Range("K99").Value = 928.44446
$ws.Range("L99").Value = 2132.7334
$ws.Range("M99").Value = 569.55554
$ws.Range("N99").Value = -5128.7334
$ws.Range("H107").Value = 1100
$ws.Range("I107").Value = 1100
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1100
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 820
$ws.Range("N107").ClearContents()
$ws.Range("H108").Value = 38000
$ws.Range("J108").Value = 38000
$ws.Range("L108").Value = 38000
$ws.Range("N108").Value = -45680

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 575.5
$ws.Range("I22").Value = 150
$ws.Range("J22").Value = 830.8
$ws.Range("K22").Value = 150
$ws.Range("L22").Value = 830.8
$ws.Range("M22").Value = 200
$ws.Range("N22").Value = -1530.8
$ws.Range("H31").Value = 2371.1304
$ws.Range("I31").Value = 1537.1621
$ws.Range("K31").Value = 1537.1621
$ws.Range("M31").Value = -1242.1621
$ws.Range("H34").Value = 2371.1304
$ws.Range("I34").Value = 1537.1621
$ws.Range("K34").Value = 1537.1621
$ws.Range("M34").Value = -1335.1621

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 721.5833
$ws.Range("J122").Value = 1999.5
$ws.Range("L122").Value = 17995.5
$ws.Range("N122").Value = -22895.5
$ws.Range("H131").Value = 3972.3823
$ws.Range("I131").Value = 443.2
$ws.Range("J131").Value = 6758.579
$ws.Range("K131").Value = 1329.6
$ws.Range("L131").Value = 20275.737
$ws.Range("M131").Value = 3710.4
$ws.Range("N131").Value = -30355.737

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5592.3887
$ws.Range("I70").Value = 4759.2593
$ws.Range("J70").Value = 8091.778
$ws.Range("K70").Value = 4759.2593
$ws.Range("L70").Value = 8091.778
$ws.Range("M70").Value = -4489.2593
$ws.Range("N70").Value = -8631.778
$ws.Range("H73").Value = 5592.3887
$ws.Range("I73").Value = 4759.2593
$ws.Range("J73").Value = 8091.778
$ws.Range("K73").Value = 4759.2593
$ws.Range("L73").Value = 8091.778
$ws.Range("M73").Value = -3823.2593
$ws.Range("N73").Value = -9963.778
$ws.Range("H80").Value = 5688.8887
$ws.Range("I80").Value = 5866.6665
$ws.Range("J80").Value = 5333.3335
$ws.Range("K80").Value = 5866.6665
$ws.Range("L80").Value = 5333.3335
$ws.Range("M80").Value = -4868.6665
$ws.Range("N80").Value = -7329.3335
$ws.Range("H83").Value = 5688.8887
$ws.Range("I83").Value = 5866.6665
$ws.Range("J83").Value = 5333.3335
$ws.Range("K83").Value = 29333.3325
$ws.Range("L83").Value = 26666.6675
$ws.Range("M83").Value = -24341.3325
$ws.Range("N83").Value = -36650.6675
$ws.Range("H97").Value = 1261.5385
$ws.Range("I97").Value = 1506.9
$ws.Range("J97").Value = 443.66666
$ws.Range("K97").Value = 1506.9
$ws.Range("L97").Value = 443.66666
$ws.Range("M97").Value = -1010.9
$ws.Range("N97").Value = -1435.66666
$ws.Range("H126").Value = 1975.8235
$ws.Range("I126").Value = 1614.5217
$ws.Range("J126").Value = 2731.2727
$ws.Range("K126").Value = 4843.5651
$ws.Range("L126").Value = 8193.8181
$ws.Range("M126").Value = -2373.5651
$ws.Range("N126").Value = -13133.8181

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1635.7906
$ws.Range("I136").Value = 1221.7354
$ws.Range("J136").Value = 3200
$ws.Range("K136").Value = 3665.2062
$ws.Range("L136").Value = 9600
$ws.Range("M136").Value = -1115.2062
$ws.Range("N136").Value = -14700

